$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# VendorCreditHeader: update the VC Date / Order Number values for row 2
# and switch those two cells (plus the header cell C1) to Text number format.
$ws1.Range("C1:C2").NumberFormat = "@"
$ws1.Range("C2").Value = "24-12-2025"
$ws1.Range("D2").Value = "231225"

# The now-unused formatted cells further down column C are removed entirely.
$ws1.Range("C3:C26").Clear()

# Update sheet view selections, and make VendorCreditItems the active tab.
$ws1.Activate()
$ws1.Range("D9").Select()
$ws2.Activate()
$ws2.Range("C19").Select()
